# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de) the two source rows
# (row 2 = the 62258c7b... file, row 3 = the e63996ec... file) move from
# "handed off" to "handed back":
#   - Status (column B, mirrored on the Overview sheet) becomes
#     "Handed back: in sync with en-US"
#   - The "Latest Target File" (E) and "Latest Handback File" (F) columns
#     are now populated with links to the same handoff/source files
#     (Source File Name / Latest Handoff File)
#   - "Latest Handback DateTime" (G) is stamped with the handback time
#   - "Handoff Reason" (H) stays "Include"

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$langSheets = @(
    @{
        Name = "zh-cn"
        HandbackTime = "2016-03-11 01:07:02"
        Row2Target = "https://github.com/OpenLocalizationTest/oltest/blob/3c54f658a638118324fc3737989443671a8f38d9/e2e/62258c7b-0558-4c87-beb7-61070730b63e.md"
        Row2Handback = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b13bc4838716069d9487f30c0c6f5e989facb5d2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/62258c7b-0558-4c87-beb7-61070730b63e.15bad8ac24728e22f947940680b0f44f7810faa6.zh-cn.xlf"
        Row3Target = "https://github.com/OpenLocalizationTest/oltest/blob/3c54f658a638118324fc3737989443671a8f38d9/e2e/e63996ec-fccc-4456-992b-59c509c21117.md"
        Row3Handback = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b13bc4838716069d9487f30c0c6f5e989facb5d2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e63996ec-fccc-4456-992b-59c509c21117.cec3c19b6d3d2ad833e939569a194b71b6d24264.zh-cn.xlf"
    },
    @{
        Name = "de-de"
        HandbackTime = "2016-03-11 01:07:23"
        Row2Target = "https://github.com/OpenLocalizationTest/oltest/blob/3c54f658a638118324fc3737989443671a8f38d9/e2e/62258c7b-0558-4c87-beb7-61070730b63e.md"
        Row2Handback = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ee85f46ef3e507b40dba02e3c6223680c25f2b3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/62258c7b-0558-4c87-beb7-61070730b63e.15bad8ac24728e22f947940680b0f44f7810faa6.de-de.xlf"
        Row3Target = "https://github.com/OpenLocalizationTest/oltest/blob/3c54f658a638118324fc3737989443671a8f38d9/e2e/e63996ec-fccc-4456-992b-59c509c21117.md"
        Row3Handback = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ee85f46ef3e507b40dba02e3c6223680c25f2b3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e63996ec-fccc-4456-992b-59c509c21117.cec3c19b6d3d2ad833e939569a194b71b6d24264.de-de.xlf"
    }
)

# Overview sheet: Status columns (B = zh-cn, C = de-de) for both file rows.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Row 2 and row 3 both flip from "Ready for handoff" to handed-back.
    $ws.Range("B2").Value = $statusHandedBack
    $ws.Range("B3").Value = $statusHandedBack

    # Pull the source file name / handoff file name so the new "Latest
    # Target File" / "Latest Handback File" links display the same file
    # names as "Source File Name" / "Latest Handoff File" (the file that
    # was handed off is the same one that is now handed back).
    $aName2 = $ws.Range("A2").Value2
    $cName2 = $ws.Range("C2").Value2
    $aName3 = $ws.Range("A3").Value2
    $cName3 = $ws.Range("C3").Value2

    $ws.Hyperlinks.Add($ws.Range("E2"), $lang.Row2Target, "", "", $aName2)
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.Row2Handback, "", "", $cName2)
    $ws.Hyperlinks.Add($ws.Range("E3"), $lang.Row3Target, "", "", $aName3)
    $ws.Hyperlinks.Add($ws.Range("F3"), $lang.Row3Handback, "", "", $cName3)

    # Stamp the handback datetime and keep the handoff reason as Include.
    $ws.Range("G2").Value = $lang.HandbackTime
    $ws.Range("G3").Value = $lang.HandbackTime
    $ws.Range("H2").Value = "Include"
    $ws.Range("H3").Value = "Include"
}
